$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row content (columns A..N) ---
$ws.Range("A1").Value = "Category"
$ws.Range("B1").Value = "Errand Type"
$ws.Range("C1").Value = "Quantity "
$ws.Range("D1").Value = "Addon1 Errand Type"
$ws.Range("E1").Value = "Addon1 Quantity"
$ws.Range("F1").Value = "Addon2 Errand Type"
$ws.Range("G1").Value = "Addon2 Quantity"
$ws.Range("H1").Value = "Terminal ID"
$ws.Range("I1").Value = "Note (Optional)"
$ws.Range("J1").Value = "Schedule?"
$ws.Range("K1").Value = "Months "
$ws.Range("L1").Value = "Day "
$ws.Range("M1").Value = "start date "
$ws.Range("N1").Value = "end date "

# --- Column widths (approximate; engine quantizes to 1/6 character units) ---
$ws.Columns.Item(1).ColumnWidth = 14.5
$ws.Columns.Item(2).ColumnWidth = 17.5
$ws.Columns.Item(3).ColumnWidth = 17.5
$ws.Columns.Item(4).ColumnWidth = 20.666666667
$ws.Columns.Item(5).ColumnWidth = 16.5
$ws.Columns.Item(6).ColumnWidth = 22
$ws.Columns.Item(7).ColumnWidth = 17.333333333
$ws.Columns.Item(8).ColumnWidth = 11.333333333
$ws.Columns.Item(9).ColumnWidth = 14.833333333
$ws.Columns.Item(10).ColumnWidth = 17.666666667
$ws.Columns.Item(11).ColumnWidth = 15.666666667
$ws.Columns.Item(12).ColumnWidth = 10.666666667
$ws.Columns.Item(13).ColumnWidth = 20
$ws.Columns.Item(14).ColumnWidth = 9.5

# --- Date number format for start/end date columns (M:N) ---
$ws.Range("M1:N1").NumberFormat = "yyyy-mm-dd;@"

# --- Font: Normal style base font Calibri -> Arial ---
$wb.Styles("Normal").Font.Name = "Arial"

# --- Selection ---
$ws.Range("C2").Select()

# --- Page setup ---
$ws.PageSetup.Orientation = 1
